$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.443.73"
$ws.Range("E2").Value = "'  -0.21%  "
$ws.Range("D3").Value = "'1.849.53"
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'240.89"
$ws.Range("E5").Value = "'  -0.93%  "
$ws.Range("D6").Value = "'0.6336"
$ws.Range("E6").Value = "'  -0.58%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D8").Value = "'3.510.25"
$ws.Range("E8").Value = "'  +89.63%  "
$ws.Range("D9").Value = "'0.07565"
$ws.Range("E9").Value = "'  +1.33%  "
$ws.Range("D10").Value = "'0.2974"
$ws.Range("E10").Value = "'  -0.53%  "
$ws.Range("D11").Value = "'24.66"
$ws.Range("E11").Value = "'  +1.56%  "
$ws.Range("D12").Value = "'3.850.71"
$ws.Range("E12").Value = "'  +84.24%  "
$ws.Range("D13").Value = "'0.07716"
$ws.Range("E13").Value = "'  +1.14%  "
$ws.Range("D14").Value = "'4.995"
$ws.Range("E14").Value = "'  -0.72%  "
$ws.Range("D15").Value = "'0.6869"
$ws.Range("E15").Value = "'  +0.19%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.000009986"
$ws.Range("E16").Value = "'  +4.71%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'83.04"
$ws.Range("E17").Value = "'  -0.72%  "
$ws.Range("D18").Value = "'6.191"
$ws.Range("E18").Value = "'  +0.52%  "
$ws.Range("D19").Value = "'29.422.17"
$ws.Range("E19").Value = "'  -0.29%  "
$ws.Range("D20").Value = "'231.87"
$ws.Range("E20").Value = "'  -1.31%  "
$ws.Range("D21").Value = "'12.51"
$ws.Range("E21").Value = "'  -0.36%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "'  +0.04%  "
$ws.Range("D23").Value = "'7.610"
$ws.Range("E23").Value = "'  -0.87%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.02%  "
$ws.Range("D25").Value = "'155.46"
$ws.Range("E25").Value = "'  -1.09%  "
$ws.Range("E26").Value = "'  -1.40%  "
$ws.Range("D27").Value = "'8.434"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("D28").Value = "'17.68"
$ws.Range("E28").Value = "'  -0.45%  "
$ws.Range("E29").Value = "'  -1.14%  "
$ws.Range("B30").Value = "RocketPoolETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D30").Value = "'3.842.63"
$ws.Range("E30").Value = "'  +91.60%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.05818"
$ws.Range("E31").Value = "'  -2.99%  "
$ws.Range("E32").Value = "'  +0.67%  "
$ws.Range("D33").Value = "'4.144"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("D34").Value = "'4.025"
$ws.Range("E34").Value = "'  -1.17%  "
$ws.Range("D35").Value = "'1.863"
$ws.Range("E35").Value = "'  -0.05%  "
$ws.Range("E36").Value = "'  -1.93%  "
$ws.Range("D37").Value = "'0.7173"
$ws.Range("E37").Value = "'  -0.23%  "
$ws.Range("D38").Value = "'2.593"
$ws.Range("E38").Value = "'  -0.10%  "
$ws.Range("D39").Value = "'1.255.11"
$ws.Range("E39").Value = "'  +4.70%  "
$ws.Range("D40").Value = "'2.795"
$ws.Range("E40").Value = "'  -0.25%  "
$ws.Range("D41").Value = "'0.01806"
$ws.Range("E41").Value = "'  +1.95%  "
$ws.Range("D42").Value = "'0.9009"
$ws.Range("E42").Value = "'  -0.81%  "
$ws.Range("D43").Value = "'6.099"
$ws.Range("E43").Value = "'  -1.01%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("E44").Value = "'  +0.05%  "
$ws.Range("D45").Value = "'101.80"
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("D46").Value = "'67.04"
$ws.Range("E46").Value = "'  +0.98%  "
$ws.Range("D47").Value = "'7.205"
$ws.Range("E47").Value = "'  -1.16%  "
$ws.Range("D48").Value = "'9.164"
$ws.Range("E48").Value = "'  +0.84%  "
$ws.Range("D49").Value = "'0.4021"
$ws.Range("E49").Value = "'  -0.33%  "
$ws.Range("D50").Value = "'1.686"
$ws.Range("E50").Value = "'  +1.88%  "
$ws.Range("E51").Value = "'  +0.19%  "
